# workshop3.pptx - "near final edits to workshops. filling in details.
# Adding comments to example code"
#
# 1. Slide 4 ("Encoding data with resistances"): fill in the title and
#    body placeholder that were left blank.
# 2. Add two new slides (5 and 6) using the same "Title and Content"
#    layout as the rest of the deck, with their title/body text filled in.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4: title + content placeholder text
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)

$title4 = $slide4.Shapes.Item(1).TextFrame.TextRange
$title4.Text = "Encoding data with resistances"
$title4.LanguageID = "en-CA"

$body4 = $slide4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "Connect push buttons in parallel to each other between GND and A0, with a different series resistance each.`r`rStart a serial connection with the Arduino, open the serial monitor`rRead from the A0 pin and print the value to the serial connection.`rPressing different buttons leads to different voltage readings on the A0 pin!`rBy cleverly selecting the resistances, you can even differentiate between presses of combinations of buttons! "
$body4.LanguageID = "en-CA"
$body4.Paragraphs(2).Text = ""
$body4.Paragraphs(4).IndentLevel = 2
$body4.Paragraphs(5).IndentLevel = 2
$body4.Paragraphs(6).IndentLevel = 2

# ---------------------------------------------------------------------
# Slide 5 (new): "Using a variable resistance to control a servo"
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Add(5, 2)

$title5 = $slide5.Shapes.Item(1).TextFrame.TextRange
$title5.Text = "Using a variable resistance to control a servo"
$title5.LanguageID = "en-CA"

$body5 = $slide5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "Use paper. Draw a line with a pencil (has to be thick and wide)`rGraphite is conductive, but has resistance. Put the GND pin at the leftmost part of your line and slide the A0 pin across the line`rResistance changes!`r`rFind the range of resistance caused by your line (serial monitor)`rWe will map the analog input to a servo position"
$body5.LanguageID = "en-CA"
$body5.Paragraphs(3).IndentLevel = 2
$body5.Paragraphs(4).IndentLevel = 2
$body5.Paragraphs(4).Text = ""
$body5.Paragraphs(6).IndentLevel = 2

# ---------------------------------------------------------------------
# Slide 6 (new): "Check out the practice script and follow along"
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Add(6, 2)

$title6 = $slide6.Shapes.Item(1).TextFrame.TextRange
$title6.Text = "Check out the practice script and follow along"
$title6.LanguageID = "en-CA"

$body6 = $slide6.Shapes.Item(2).TextFrame.TextRange
$body6.Text = "practice script " + [char]8220 + "servo-example-variable-resistance.ino" + [char]8221
$body6.LanguageID = "en-CA"
